## Azure_Network_Design.pptx - "Updated routing in Azure with service endpoint"
##
## Re-numbers the spoke VNet from 10.1.0.0/16 to 10.3.0.0/16 on the network
## diagram slide (slide 2): the subnet label, the spoke VNet label, and the
## service-endpoint-enabled VM's address label.

$p = $ppt.ActivePresentation

# --- Slide 2: network diagram -------------------------------------------
$s = $p.Slides.Item(2)

# "Group 54" holds the Application-subnet box and the Spoke VNet box.
$grp = $s.Shapes.Item(7)

$subnetLabel = $grp.GroupItems.Item(5)
if ($subnetLabel.TextFrame.TextRange.Text -eq "Subnet – 10.1.0.0/24 (Application)") {
    $subnetLabel.TextFrame.TextRange.Text = "Subnet – 10.3.0.0/24 (Application)"
}

$vnetLabel = $grp.GroupItems.Item(6)
if ($vnetLabel.TextFrame.TextRange.Text -eq "Spoke VNet – 10.1.0.0/16") {
    $vnetLabel.TextFrame.TextRange.Text = "Spoke VNet – 10.3.0.0/16"
}

# "Rectangle 106" is the small IP-address label (the VM behind the
# service endpoint) near the top-left of the spoke VNet.
$ipLabel = $s.Shapes.Item(10)
if ($ipLabel.TextFrame.TextRange.Text -eq "10.1.0.4") {
    $ipLabel.TextFrame.TextRange.Text = "10.3.0.4"
}

# --- Refresh the cached "today" date field shown on every layout/master -
# (PowerPoint recomputes the auto date placeholder whenever the deck is
# edited and saved; the source deck was re-saved two days later.)
$newDate = "2021-05-23"

$master = $p.SlideMaster
$masterDate = $master.Shapes.Item(3)
if ($masterDate.TextFrame.TextRange.Text -eq "2021-05-21") {
    $masterDate.TextFrame.TextRange.Text = $newDate
}

$dateIdxByLayout = @{1=3; 2=3; 3=3; 4=4; 5=6; 6=2; 7=1; 8=4; 9=4; 10=3; 11=3}
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $idx = $dateIdxByLayout[$i]
    $dateShape = $layout.Shapes.Item($idx)
    if ($dateShape.TextFrame.TextRange.Text -eq "2021-05-21") {
        $dateShape.TextFrame.TextRange.Text = $newDate
    }
}
